$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Wei%Xia%NULL%2,               Jianbo%Shao%xzhuge27@163.com%2,               Yu%Guo%NULL%4,               Yu%Guo%NULL%0,               Xuehua%Peng%NULL%4,               Xuehua%Peng%NULL%0,               Zhen%Li%NULL%2,               Daoyu%Hu%NULL%2]"
$ws.Range("E3").Value = "[Wei%Xia%NULL%0,               Jianbo%Shao%xzhuge27@163.com%0,               Yu%Guo%NULL%0,               Yu%Guo%NULL%0,               Xuehua%Peng%NULL%0,               Xuehua%Peng%NULL%0,               Zhen%Li%NULL%0,               Daoyu%Hu%NULL%0]"
$ws.Range("E4").Value = "[Michela%Antonelli%NULL%1,               Rose S%Penfold%NULL%1,               Jordi%Merino%NULL%1,               Carole H%Sudre%NULL%1,               Erika%Molteni%NULL%1,               Sarah%Berry%NULL%1,               Liane S%Canas%NULL%1,               Mark S%Graham%NULL%1,               Kerstin%Klaser%NULL%1,               Marc%Modat%NULL%1,               Benjamin%Murray%NULL%1,               Eric%Kerfoot%NULL%1,               Liyuan%Chen%NULL%1,               Jie%Deng%NULL%1,               Marc F%Österdahl%NULL%1,               Nathan J%Cheetham%NULL%1,               David A%Drew%NULL%1,               Long H%Nguyen%NULL%1,               Joan Capdevila%Pujol%NULL%1,               Christina%Hu%NULL%1,               Somesh%Selvachandran%NULL%1,               Lorenzo%Polidori%NULL%1,               Anna%May%NULL%1,               Jonathan%Wolf%NULL%1,               Andrew T%Chan%NULL%1,               Alexander%Hammers%NULL%1,               Emma L%Duncan%NULL%1,               Tim D%Spector%NULL%1,               Sebastien%Ourselin%NULL%1,               Claire J%Steves%claire.j.steves@kcl.ac.uk%1]"
$ws.Range("E5").Value = "[Qinxue%Shen%NULL%1,               Wei%Guo%NULL%1,               Ting%Guo%NULL%1,               Jinhua%Li%NULL%1,               Wenlong%He%NULL%1,               Shanshan%Ni%NULL%1,               Xiaoli%Ouyang%NULL%1,               Jiyang%Liu%NULL%1,               Yuanlin%Xie%NULL%1,               Xin%Tan%NULL%1,               Zhiguo%Zhou%13807311490@163.com%1,               Hong%Peng%penghong66@csu.edu.cn%1]"
$ws.Range("E6").Value = "[Ya\u2010nan%Han%xref no email%1,           Zhan\u2010wei%Feng%xref no email%1,           Li\u2010na%Sun%xref no email%1,           Xiao\u2010xia%Ren%xref no email%1,           Hua%Wang%xref no email%0,           Yong\u2010ming%Xue%xref no email%1,           Yi%Wang%xref no email%1,           Ying%Fang%xref no email%1]"
$ws.Range("E7").Value = "[Lan%Zhang%NULL%1,               Songming%Huang%NULL%1]"
$ws.Range("E8").Value = "[Lu%Lin%NULL%1,               Xiayang%Jiang%NULL%1,               Zhenling%Zhang%NULL%1,               Siwen%Huang%NULL%1,               Zhenyi%Zhang%NULL%1,               Zhaoxiong%Fang%NULL%1,               Zhiqiang%Gu%NULL%1,               Liangqing%Gao%NULL%1,               Honggang%Shi%NULL%1,               Lei%Mai%NULL%1,               Yuan%Liu%NULL%1,               Xianqi%Lin%NULL%1,               Renxu%Lai%NULL%1,               Zhixiang%Yan%NULL%1,               Xiaofeng%Li%NULL%2,               Xiaofeng%Li%NULL%0,               Hong%Shan%NULL%1]"
$ws.Range("E9").Value = "[Dan%Sun%NULL%0,               Hui%Li%NULL%1,               Xiao-Xia%Lu%NULL%1,               Han%Xiao%NULL%1,               Jie%Ren%NULL%1,               Fu-Rong%Zhang%792523496@qq.com%1,               Zhi-Sheng%Liu%liuzsc@126.com%1]"
$ws.Range("E10").Value = "[Yi%Xu%NULL%1,               Xufang%Li%NULL%1,               Bing%Zhu%NULL%1,               Huiying%Liang%NULL%1,               Chunxiao%Fang%NULL%1,               Yu%Gong%NULL%1,               Qiaozhi%Guo%NULL%1,               Xin%Sun%NULL%1,               Danyang%Zhao%NULL%1,               Jun%Shen%NULL%1,               Huayan%Zhang%NULL%1,               Hongsheng%Liu%NULL%2,               Hongsheng%Liu%NULL%0,               Huimin%Xia%huiminxia@hotmail.com%1,               Jinling%Tang%jltang@cuhk.edu.hk%1,               Kang%Zhang%kang.zhang@gmail.com%1,               Sitang%Gong%sitangg@126.com%2,               Sitang%Gong%sitangg@126.com%0]"
$ws.Range("E11").Value = "[Alfredo%Tagarro%xref no email%1,           Cristina%Epalza%xref no email%1,           Mar%Santos%xref no email%1,           Francisco Jos\u00e9%Sanz-Santaeufemia%xref no email%1,           Enrique%Otheo%xref no email%1,           Cinta%Moraleda%xref no email%1,           Cristina%Calvo%xref no email%1]"
$ws.Range("E14").Value = "[Yuhan%Xing%NULL%1,               Wei%Ni%NULL%1,               Qin%Wu%NULL%1,               Wenjie%Li%NULL%1,               Guoju%Li%NULL%1,               Wendi%Wang%NULL%1,               Jianning%Tong%NULL%1,               Xiufeng%Song%NULL%1,               Gary Wing Kin%Wong%NULL%1,               Quansheng%Xing%NULL%1]"
$ws.Range("E15").Value = "[Giacomo%Brisca%giacomobrisca@gmail.com%1,               Marta%Ferretti%NULL%2,               Marta%Ferretti%NULL%0,               Giulia%Sartoris%NULL%1,               Maria Beatrice%Damasio%NULL%1,               Isabella%Buffoni%NULL%2,               Isabella%Buffoni%NULL%0,               Daniela%Pirlo%NULL%1,               Marta%Romanengo%NULL%2,               Emanuela%Piccotti%NULL%1]"
$ws.Range("E16").Value = "[Louise%Gaborieau%NULL%1,               Celine%Delestrain%NULL%1,               Philippe%Bensaid%NULL%1,               Audrey%Vizeneux%NULL%1,               Philippe%Blanc%NULL%1,               Aurélie%Garraffo%NULL%1,               Emilie%Georget%NULL%1,               Arnaud%Chalvon%NULL%1,               Nathalie%Garrec%NULL%1,               Yacine%Laoudi%NULL%1,               Emmanuelle%Varon%NULL%1,               Sébastien%Rouget%NULL%1,               Alexandre%Pupin%NULL%1,               Khaled%Abdel Aal%NULL%1,               David%Toulorge%NULL%1,               Sarah%Ducrocq%NULL%1,               Catherine%Barrey%NULL%1,               Letitia%Pantalone%NULL%1,               Blandine%Robert%NULL%1,               Lydie%Joly-Sanchez%NULL%1,               Caroline%Thach%NULL%2,               Caroline%Thach%NULL%0,               Caroline%Masserot-Lureau%NULL%1,               Jamilé%Chahine%NULL%1,               Veronica%Risso Garcia-Roudaut%NULL%1,               Jonathan%Rozental%NULL%1,               Sylvie%Nathanson%NULL%1,               Mohamed%Khaled%NULL%1,               Alexis%Mandelcwajg%NULL%1,               Nadia%Demayer%NULL%1,               Stéphanie%Muller%NULL%1,               Mustapha%Mazerghane%NULL%1,               Ralph%Epaud%NULL%1,               Béatrice%Pellegrino%NULL%2,               Béatrice%Pellegrino%NULL%0,               Fouad%Madhi%NULL%1]"
$ws.Range("E17").Value = "[Chun‐Zhen%Hua%NULL%1,               Zi‐Ping%Miao%NULL%1,               Ji‐Shan%Zheng%NULL%2,               Ji‐Shan%Zheng%NULL%0,               Qian%Huang%NULL%1,               Qing‐Feng%Sun%NULL%1,               Hong‐Ping%Lu%NULL%1,               Fei‐Fei%Su%NULL%1,               Wei‐Hong%Wang%NULL%1,               Lie‐Ping%Huang%NULL%1,               Da‐Qing%Chen%NULL%1,               Zhi‐Wei%Xu%NULL%1,               Le‐Dan%Ji%NULL%1,               Hong‐Ping%Zhang%NULL%1,               Xiao‐Wei%Yang%NULL%1,               Ming‐Hui%Li%NULL%1,               Yue‐Yan%Mao%NULL%1,               Man‐Zhen%Ying%NULL%1,               Sheng%Ye%NULL%1,               Qiang%Shu%shuqiang@zju.edu.cn%1,               En‐Fu%Chen%enfchen@cdc.zj.cn%2,               En‐Fu%Chen%enfchen@cdc.zj.cn%0,               Jian‐Feng%Liang%NULL%1,               Wei%Wang%NULL%1,               Zhi‐Min%Chen%NULL%1,               Wei%Li%NULL%1,               Jun‐Fen%Fu%fjf68@zju.edu.cn%1]"
$ws.Range("E18").Value = "[María%de Ceano-Vivas%NULL%1,               Irene%Martín-Espín%NULL%1,               Teresa%del Rosal%NULL%1,               Marta%Bueno-Barriocanal%NULL%2,               Marta%Bueno-Barriocanal%NULL%0,               Marta%Plata-Gallardo%NULL%1,               José Antonio%Ruiz-Domínguez%NULL%1,               Rosario%López-López%NULL%2,               Rosario%López-López%NULL%0,               Miguel Ángel%Molina-Gutiérrez%NULL%1,               Patricia%Bote-Gascón%NULL%2,               Patricia%Bote-Gascón%NULL%0,               Isabel%González-Bertolín%NULL%1,               Paula%García-Sánchez%NULL%2,               Paula%García-Sánchez%NULL%0,               Julia%Martín-Sánchez%NULL%1,               Begoña%de Miguel-Lavisier%NULL%1,               Talía%Sainz%NULL%1,               Fernando%Baquero-Artigao%NULL%1,               Ana%Méndez-Echevarría%NULL%1,               Cristina%Calvo%NULL%2,               Cristina%Calvo%NULL%0]"
$ws.Range("E19").Value = "[Niccolò%Parri%niccolo.parri@meyer.it%1,              Anna Maria%Magistà%annamaria.magista@auslromagna.it%1,              Federico%Marchetti%federico.marchetti@auslromagna.it%1,              Barbara%Cantoni%barbara.cantoni@policlinico.mi.it%1,              Alberto%Arrighini%alberto.arrighini@asst-spedalicivili.it%1,              Marta%Romanengo%marta.romanengo@gmail.com%0,              Enrico%Felici%enrico.felici@ospedale.al.it%1,              Antonio%Urbino%afurbino@icloud.com%1,              Liviana%Da Dalt%liviana.dadalt@unipd.it%1,              Lucio%Verdoni%lverdoni@asst-pg23.it%1,              Benedetta%Armocida%benedetta.armocida@burlo.trieste.it%1,              Benedetta%Covi%benedetta.covi@burlo.trieste.it%1,              Ilaria%Mariani%ilaria.mariani@burlo.trieste.it%1,              Roberta%Giacchero%Roberta.Giacchero@asst-lodi.it%1,              Anna Maria%Musolino%annamaria.musolino@opbg.net%1,              Marco%Binotti%marco.binotti@med.uniupo.it%1,              Paolo%Biban%paolo.biban@aovr.veneto.it%1,              Silvia%Fasoli%silvia.fasoli@asst-mantova.it%1,              Chiara%Pilotto%chiara.pilotto@asufc.sanita.fvg.it%1,              Flavia%Nicoloso%flavianicoloso@gmail.com%1,              Massimiliano%Raggi%massimiali.raggi@apss.tn.it%1,              Elisabetta%Miorin%elisabetta.miorin@asufc.sanita.fvg.it%1,              Danilo%Buonsenso%danilobuonsenso@gmail.com%1,              Massimo%Chiossi%max.chiossi@gmail.com%1,              Rino%Agostiniani%rinoagostiniani@gmail.com%1,              Anna%Plebani%annamaria.plebani@asst-settelaghi.it%1,              Maria Antonietta%Barbieri%mantonietta.barbieri@opbg.net%1,              Marcello%Lanari%marcello.lanari@unibo.it%1,              Serena%Arrigo%arrigoserena@gmail.com%1,              Elena%Zoia%elena.zoia@asst.fbf.sacco.it%1,              Matteo%Lenge%matteo.lenge@meyer.it%1,              Stefano%Masi%stefano.masi@meyer.it%1,              Egidio%Barbi%egidio.barbi@burlo.trieste.it%1,              Marzia%Lazzerini%marzia.lazzerini@burlo.trieste.it%1,              NULL%NULL%NULL%1]"
